$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00348977135980746
$ws.Range("C2").Value = 0.00601684717208183
$ws.Range("D2").Value = 0.990914560770156
$ws.Range("E2").Value = 0.00240673886883273
$ws.Range("F2").Value = 0.989891696750902
$ws.Range("G2").Value = 0.00667870036101083
$ws.Range("I2").Value = 0.00246690734055355
$ws.Range("J2").Value = 0.0000601684717208183
$ws.Range("K2").Value = 0.0154632972322503
$ws.Range("L2").Value = 0.648134777376655
$ws.Range("M2").Value = 0.00276774969915764
$ws.Range("N2").Value = 0.00469314079422383
$ws.Range("O2").Value = 0.998916967509025
$ws.Range("P2").Value = 0.000300842358604091
$ws.Range("R2").Value = 0.999879663056558
$ws.Range("S2").Value = 0.00246690734055355
$ws.Range("T2").Value = 0.000421179302045728
$ws.Range("U2").Value = 0.00547533092659446
$ws.Range("W2").Value = 0.0755716004813478
$ws.Range("X2").Value = 0.00716004813477738
$ws.Range("C3").Value = 0.000120336943441637
$ws.Range("D3").Value = 0.00601684717208183
$ws.Range("E3").Value = 0.033453670276775
$ws.Range("F3").Value = 0.00607701564380265
$ws.Range("H3").Value = 0.648977135980746
$ws.Range("I3").Value = 0.990673886883273
$ws.Range("J3").Value = 0.991335740072202
$ws.Range("K3").Value = 0.051323706377858
$ws.Range("L3").Value = 0.000300842358604091
$ws.Range("M3").Value = 0.628880866425993
$ws.Range("N3").Value = 0.926113116726835
$ws.Range("P3").Value = 0.000120336943441637
$ws.Range("Q3").Value = 0.961672683513839
$ws.Range("S3").Value = 0.000481347773766546
$ws.Range("T3").Value = 0.0123947051744886
$ws.Range("U3").Value = 0.799157641395909
$ws.Range("V3").Value = 0.00276774969915764
$ws.Range("W3").Value = 0.00613718411552347
$ws.Range("X3").Value = 0.00246690734055355
$ws.Range("B4").Value = 0.990493381468111
$ws.Range("C4").Value = 0.00300842358604091
$ws.Range("D4").Value = 0.00288808664259928
$ws.Range("E4").Value = 0.000120336943441637
$ws.Range("F4").Value = 0.000300842358604091
$ws.Range("G4").Value = 0.993321299638989
$ws.Range("H4").Value = 0.00246690734055355
$ws.Range("I4").Value = 0.00601684717208183
$ws.Range("K4").Value = 0.00409145607701564
$ws.Range("L4").Value = 0.325150421179302
$ws.Range("M4").Value = 0.0036101083032491
$ws.Range("N4").Value = 0.0000601684717208183
$ws.Range("O4").Value = 0.000782190132370638
$ws.Range("P4").Value = 0.999518652226233
$ws.Range("R4").Value = 0.000120336943441637
$ws.Range("U4").Value = 0.00120336943441637
$ws.Range("V4").Value = 0.000782190132370638
$ws.Range("W4").Value = 0.913056558363418
$ws.Range("X4").Value = 0.990373044524669
$ws.Range("B5").Value = 0.00601684717208183
$ws.Range("C5").Value = 0.990854392298436
$ws.Range("D5").Value = 0.000120336943441637
$ws.Range("E5").Value = 0.962755716004813
$ws.Range("F5").Value = 0.00330926594464501
$ws.Range("H5").Value = 0.341395908543923
$ws.Range("I5").Value = 0.000842358604091456
$ws.Range("J5").Value = 0.00860409145607702
$ws.Range("K5").Value = 0.925210589651023
$ws.Range("L5").Value = 0.0242478941034898
$ws.Range("M5").Value = 0.361672683513839
$ws.Range("N5").Value = 0.0659446450060168
$ws.Range("O5").Value = 0.000180505415162455
$ws.Range("Q5").Value = 0.0365824308062575
$ws.Range("S5").Value = 0.99705174488568
$ws.Range("T5").Value = 0.985619735258724
$ws.Range("U5").Value = 0.188688327316486
$ws.Range("V5").Value = 0.996450060168472
$ws.Range("W5").Value = 0.00192539109506619
